$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 2
$ws.Range("F3").Value = 10982
$ws.Range("F4").Value = 265
$ws.Range("F5").Value = 1221
$ws.Range("F6").Value = 1094
$ws.Range("F7").Value = 848
$ws.Range("F8").Value = 287
$ws.Range("F10").Value = 1174
$ws.Range("F13").Value = 896
$ws.Range("F15").Value = 2041
$ws.Range("F17").Value = 1008
$ws.Range("F18").Value = 839
$ws.Range("F19").Value = 559
$ws.Range("F20").Value = 809
$ws.Range("F21").Value = 925
$ws.Range("F24").Value = 640
$ws.Range("F25").Value = 660
$ws.Range("F26").Value = 128
$ws.Range("F27").Value = 357
$ws.Range("F28").Value = 1020
$ws.Range("F29").Value = 48
$ws.Range("F30").Value = 501
$ws.Range("F31").Value = 178
$ws.Range("F32").Value = 255
$ws.Range("F33").Value = 243
$ws.Range("F34").Value = 586
$ws.Range("F35").Value = 1883
$ws.Range("F36").Value = 394
$ws.Range("F37").Value = 37
$ws.Range("F38").Value = 1448
$ws.Range("F41").Value = 51
$ws.Range("F45").Value = 78
$ws.Range("F49").Value = 83

$ws = $wb.Worksheets.Item(2)
$ws.Range("F4").Value = 90
$ws.Range("F5").Value = 199
$ws.Range("F12").Value = 171
$ws.Range("F14").Value = 140

$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 2183
$ws.Range("F3").Value = 638
$ws.Range("F4").Value = 577

$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 2183
$ws.Range("F3").Value = 638
$ws.Range("F5").Value = 10982
$ws.Range("F6").Value = 265
$ws.Range("F7").Value = 90
$ws.Range("F8").Value = 577
$ws.Range("F9").Value = 1094
$ws.Range("F10").Value = 199
$ws.Range("F11").Value = 1174
$ws.Range("F14").Value = 896
$ws.Range("F15").Value = 2041
$ws.Range("F17").Value = 1008
$ws.Range("F18").Value = 839
$ws.Range("F19").Value = 559
$ws.Range("F20").Value = 809
$ws.Range("F21").Value = 925
$ws.Range("F25").Value = 640
$ws.Range("F28").Value = 660
$ws.Range("F29").Value = 128
$ws.Range("F30").Value = 357
$ws.Range("F31").Value = 1020
$ws.Range("F33").Value = 48
$ws.Range("F34").Value = 501
$ws.Range("F35").Value = 178
$ws.Range("F36").Value = 255
$ws.Range("F37").Value = 243
$ws.Range("F38").Value = 37
$ws.Range("F39").Value = 1448
$ws.Range("F42").Value = 51
$ws.Range("F46").Value = 78
$ws.Range("F48").Value = 83

$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("C17").Value = "广州·音漫派国湿演唱会-《狐妖小红娘》《一人之下》领衔国漫原声音乐现场"